# OLX Monitor 2026-02-17 13:19
# Appends a fresh monitoring snapshot (same 8 listings, new check timestamp)
# to the bottom of the "PODSUMOWANIE" detail log, mirroring the most recent
# prior batch (rows 23:30) down into rows 31:38.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Duplicate the last logged batch (values + styles) into the next 8 rows.
$src = $ws.Range("A23:H30")
$dst = $ws.Range("A31")
$src.Copy($dst)

# Stamp the new batch with this run's "last checked" timestamp.
$ws.Range("A31:A38").Value = "2026-02-17 13:19:00"

# The "Ładny pokój jednoosobowy..." listing (dawnypatron) aged one more day.
$ws.Range("F37").Value = 515
